$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.780.88"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "2.028.84"
$ws.Range("E3").Value = "  -1.70%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.608"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.14"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.99%  "
$ws.Range("E9").Value = "  -0.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0813"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.56"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.79%  "
$ws.Range("D13").Value = "2.328.57"
$ws.Range("E13").Value = "  -1.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("E16").Value = "  -1.87%  "
$ws.Range("D17").Value = "2.033.10"
$ws.Range("E17").Value = "  -1.83%  "
$ws.Range("D18").Value = "37.697.56"
$ws.Range("E18").Value = "  -0.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("D21").Value = "0.0₃0824"
$ws.Range("E21").Value = "  -1.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.93%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").Value = "  -2.67%  "
$ws.Range("E25").Value = "  -1.89%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.79%  "
$ws.Range("E28").Value = "  -3.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.92%  "
$ws.Range("E30").Value = "  -5.02%  "
$ws.Range("E31").Value = "  +0.85%  "
$ws.Range("E32").Value = "  -2.98%  "
$ws.Range("E33").Value = "  +5.33%  "
$ws.Range("E34").Value = "  -3.00%  "
$ws.Range("E35").Value = "  -2.36%  "
$ws.Range("E36").Value = "  +5.13%  "
$ws.Range("E37").Value = "  -4.46%  "
$ws.Range("E38").Value = "  -3.04%  "
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").Value = "1.530.88"
$ws.Range("E40").Value = "  +3.38%  "
$ws.Range("E41").Value = "  -1.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.62%  "
$ws.Range("B45").Value = "HuobiToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.51%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0916"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.15%  "
$ws.Range("E47").Value = "  -2.00%  "
$ws.Range("E48").Value = "  -1.69%  "
$ws.Range("E49").Value = "  -0.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").Value = "2.218.39"
$ws.Range("E51").Value = "  -1.62%  "
